$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Paragraph: "Created an operating system class ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "Created an operating system class that controlled time through a scheduler. This value was passed to each algorithm."
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, "This group c", 2)

$r.InsertXML("<w:p $w><w:r><w:t>reated an operating system class that controlled time through a</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`">n </w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:proofErr w:type=`"spellStart`"/><w:r><w:t>enum</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> called</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> schedu</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>ler. This enumerated value is then passed into each algorithm for execution</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>.</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> This value is the incrementation value as each process goes through a loop. This incrementation simulates the milliseconds of time that a process executes.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "Created a process object which is sent through the scheduler ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "Created a process object which is sent through the scheduler to control CPU ticks and execution."
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, "They c", 2)

$r.InsertXML("<w:p $w><w:r><w:t>reated a process object which is sent through the scheduler to control CPU ticks and execution.</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> Processes are modeled as a class object that has fields such as ID, bursts, and switch times. This allows the calculation of each statistic very carefully based on what information is needed.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "Each algorithm is its own class and it inherits from the scheduler. ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "Each algorithm is its own class and it inherits from the scheduler. That way time can be controlled for execution."
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, "Each algorithm is its own class and it inheri", 2)

$r.InsertXML("<w:p $w><w:r><w:t>ts from the scheduler. This allows each algorithm to utilize the same definition and variables associated with time. It also in a way simulates each algorithm running on the same machine. Also, each algorithm does not need to have time defined within.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "They used one function to run processes ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "They used one function to run processes and used nested loops with a higher counter to simulate multiple core FCFS."
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, $orig, 2)

$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:proofErr w:type=`"gramStart`"/><w:r><w:t>So</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> they could re-use code for the FCFS algorithm and </w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>adapt it to a multiprocessor simulation with a higher counter in the loop such as 8 for 4 cores.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "One of the major issues occurred in FCFS ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "One of the major issues occurred in FCFS when one process finished and had a very short I/O burst, it was put back into the queue, but was also still in the execution vector. Their solution was to remove it from the vector when it got close to finishing not when it was finished. Also an issue was adapting the single runProcess function to all different algorithms and multiprocessor approaches."
$chunk1 = "One of the major issues occurred in FCFS when one process finished and had a very short I/O burst, it was put back into the queue, but was also still in the execution vector. Their solution "
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, $chunk1, 2)

$r.InsertXML("<w:p $w><w:r><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">was to remove it from the vector when it got close to finishing not when it was finished. </w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:proofErr w:type=`"gramStart`"/><w:r><w:t>Also</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> an issue was adapting the single </w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:proofErr w:type=`"spellStart`"/><w:r><w:t>runProcess</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> function to all different algorithms </w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>and multiprocessor approaches. Also, there was no destructor for the scheduler so that was addressed. Because of this there was a small memory leak.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "Arrays, Vectors were used for execution. ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "Arrays, Vectors were used for execution. The Scheduler used a queue and a map was used for a process table. It would map to the pointer for each process."
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, "A", 2)

$r.InsertXML("<w:p $w><w:r><w:t>rrays and</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> Vectors were used for execution. The Scheduler used a queue and a map was used fo</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>r a process table. Each process is stored via a pointer to a location within this map.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "A stack was considered to sort processes differently. ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "A stack was considered to sort processes differently. But first in first out was needed so a queue was needed."
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, "A stack was considered", 2)

$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> to sort processes differently b</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>ut</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>,</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> first in </w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>first sorting was required so a queue was chosen to be implemented as well as a map</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t>.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "If utilizing the file system a heap may have been applicable. "
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "If utilizing the file system a heap may have been applicable. "
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, "If utilizing the file ", 2)

$r.InsertXML("<w:p $w><w:r><w:t>system,</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> a heap may have been applicable. </w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "Followed class hierarchy well from the UML. ..." (drop leading
# lastRenderedPageBreak, expand closing sentence)
# ---------------------------------------------------------------------------
$r = $d.Content
$orig = "Followed class hierarchy well from the UML. Everything works from the Scheduler class. "
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, $orig, 2)

$r = $d.Content
$orig2 = [char]0x2019
$orig2text = "One thing is that FCFS and SPN look very similar so they could probably be merged. Multicore" + $orig2 + "s number of cores is only used in the FCFS algorithm so this could be added to all other algorithms."
$chunk1 = "One thing is that FCFS and SPN look very similar so they could probably be merged. Multicore" + $orig2 + "s number of cores is only used in the FCFS "
$null = $r.Find.Execute($orig2text, $false, $false, $false, $false, $false, $true, 1, $false, $chunk1, 2)

$r.InsertXML("<w:p $w><w:r><w:t>algorithm so this could be expanded</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> to</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> be applicable for</w:t></w:r></w:p>")
$r.InsertXML("<w:p $w><w:r><w:t xml:space=`"preserve`"> all other algorithms.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# Paragraph: "Exceptions are handled through a try, catch block in main. ..."
# Move the existing _GoBack bookmark from the end of the paragraph to the
# middle, splitting the run.
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$r = $d.Content
$orig = "Exceptions are handled through a try, catch block in main. Errors were fixed when they arose as opposed to try catch."
$null = $r.Find.Execute($orig, $false, $false, $false, $false, $false, $true, 1, $false, "Exceptions ar", 2)

$r.InsertXML("<w:p $w><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t>e handled through a try, catch block in main. Errors were fixed when they arose as opposed to try catch.</w:t></w:r></w:p>")
